$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F column "想去人数" (number of attendees interested)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 64
$wsExhibit.Range("F4").Value = 162
$wsExhibit.Range("F6").Value = 5369
$wsExhibit.Range("F8").Value = 5380
$wsExhibit.Range("F9").Value = 628
$wsExhibit.Range("F11").Value = 1374
$wsExhibit.Range("F12").Value = 6
$wsExhibit.Range("F13").Value = 108

# Sheet "全部类型" (sheet4): same underlying records, update matching F column values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 64
$wsAll.Range("F4").Value = 162
$wsAll.Range("F7").Value = 5369
$wsAll.Range("F9").Value = 5380
$wsAll.Range("F10").Value = 628
$wsAll.Range("F12").Value = 1374
$wsAll.Range("F13").Value = 6
$wsAll.Range("F14").Value = 108
